# The species records stored in rows 43/45, 66/67 and 78/79 of the
# "Artfynd" sheet got reordered. All of the record-specific data
# (everything except the shared location/observer columns, which are
# identical between the paired rows anyway) needs to be swapped
# between each pair of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns holding per-record data that differs between the two rows
# of each swapped pair.
$columns = @("A","B","D","E","F","G","H","I","J","Q","R","S","Z","AB","AC","AJ","AK","AO")

# Column "I" ("Antal") stores numeric-looking values as text in this
# workbook, so it needs to be forced to text so Excel does not
# silently reinterpret "50" as the number 50.
$textColumns = @("I")

function Swap-Rows($ws, $row1, $row2, $columns, $textColumns) {
    # Read all the old values first so overwriting one side does not
    # affect reading the other.
    $values1 = @{}
    $values2 = @{}
    foreach ($col in $columns) {
        $values1[$col] = $ws.Range("$col$row1").Value2
        $values2[$col] = $ws.Range("$col$row2").Value2
    }

    foreach ($col in $columns) {
        $cell1 = $ws.Range("$col$row1")
        $cell2 = $ws.Range("$col$row2")
        $new1 = $values2[$col]
        $new2 = $values1[$col]

        if (($textColumns -contains $col) -and (($new1) -or ($new2))) {
            # Force text storage so numeric-looking values (e.g. "50")
            # are not silently reinterpreted as numbers.
            $cell1.NumberFormat = "@"
            $cell2.NumberFormat = "@"
            $cell1.Value2 = $new1
            $cell2.Value2 = $new2
            $cell1.Style = "Normal"
            $cell2.Style = "Normal"
        } else {
            $cell1.Value2 = $new1
            $cell2.Value2 = $new2
        }
    }
}

Swap-Rows $ws 43 45 $columns $textColumns
Swap-Rows $ws 66 67 $columns $textColumns
Swap-Rows $ws 78 79 $columns $textColumns
